$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Key / (blank) / Type changes from varchar -> integer
$ws.Range("C2").Value = "integer"

# Row 4: rare_diseases_family_id / Rare Diseases Family Id -> referral_id / Referral_id
$ws.Range("A4").Value = "referral_id"
$ws.Range("B4").Value = "Referral_id"

# Row 7: sample_id / Sample Id -> platekey / platekey
$ws.Range("A7").Value = "platekey"
$ws.Range("B7").Value = "platekey"
